$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H137").Value = 1188318.8
$ws.Range("I137").Value = 2141451
$ws.Range("J137").Value = 5120.207
$ws.Range("K137").Value = 6424353
$ws.Range("L137").Value = 15360.621
$ws.Range("M137").Value = -6421803
$ws.Range("N137").Value = -20460.621
$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16457.227
$ws.Range("I32").Value = 15387.436
$ws.Range("J32").Value = 24801.6
$ws.Range("K32").Value = 15387.436
$ws.Range("L32").Value = 24801.6
$ws.Range("M32").Value = -15100.436
$ws.Range("N32").Value = -25375.6
$ws.Range("H45").Value = 1846.4595
$ws.Range("I45").Value = 1725.1786
$ws.Range("J45").Value = 2223.7778
$ws.Range("K45").Value = 1725.1786
$ws.Range("L45").Value = 2223.7778
$ws.Range("M45").Value = -1348.1786
$ws.Range("N45").Value = -2977.7778
$ws.Range("H97").Value = 1903.4546
$ws.Range("I97").Value = 1770.8889
$ws.Range("J97").Value = 2500
$ws.Range("K97").Value = 1770.8889
$ws.Range("L97").Value = 2500
$ws.Range("M97").Value = -1274.8889
$ws.Range("N97").Value = -3492

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 41933.332
$ws.Range("J59").Value = 41933.332
$ws.Range("L59").Value = 41933.332
$ws.Range("N59").Value = -43627.332
$ws.Range("H105").Value = 2324.16
$ws.Range("I105").Value = 2177.4783
$ws.Range("J105").Value = 4011
$ws.Range("K105").Value = 2177.4783
$ws.Range("L105").Value = 4011
$ws.Range("M105").Value = -430.4783000000002
$ws.Range("N105").Value = -7505
$ws.Range("H107").Value = 2161.1143
$ws.Range("I107").Value = 1850.8667
$ws.Range("J107").Value = 4022.6
$ws.Range("K107").Value = 1850.8667
$ws.Range("L107").Value = 4022.6
$ws.Range("M107").Value = 69.13329999999996
$ws.Range("N107").Value = -7862.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3678.3438
$ws.Range("I5").Value = 12850.875
$ws.Range("J5").Value = 620.8333
$ws.Range("K5").Value = 38552.625
$ws.Range("L5").Value = 1862.4999
$ws.Range("M5").Value = -38440.625
$ws.Range("N5").Value = -2086.4999
$ws.Range("H68").Value = 1226.4247
$ws.Range("I68").Value = 963.8889
$ws.Range("J68").Value = 1312.3455
$ws.Range("K68").Value = 2891.6667
$ws.Range("L68").Value = 3937.0365
$ws.Range("M68").Value = -2080.6667
$ws.Range("N68").Value = -5559.0365
$ws.Range("H71").Value = 1226.4247
$ws.Range("I71").Value = 963.8889
$ws.Range("J71").Value = 1312.3455
$ws.Range("K71").Value = 8675.000100000001
$ws.Range("L71").Value = 11811.1095
$ws.Range("M71").Value = -4619.000100000001
$ws.Range("N71").Value = -19923.1095
$ws.Range("H80").Value = 34571024
$ws.Range("J80").Value = 35734276
$ws.Range("L80").Value = 107202828
$ws.Range("N80").Value = -107204700
$ws.Range("H83").Value = 34571024
$ws.Range("J83").Value = 35734276
$ws.Range("L83").Value = 321608484
$ws.Range("N83").Value = -321617844
$ws.Range("H97").Value = 101350
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 101350
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 304050
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -305042
$ws.Range("H113").Value = 4451.222
$ws.Range("I113").Value = 7705.643
$ws.Range("J113").Value = 946.46155
$ws.Range("K113").Value = 23116.929
$ws.Range("L113").Value = 2839.38465
$ws.Range("M113").Value = -20946.929
$ws.Range("N113").Value = -7179.38465
$ws.Range("H114").Value = 1592.3846
$ws.Range("I114").Value = 1205.75
$ws.Range("J114").Value = 1923.7858
$ws.Range("K114").Value = 3617.25
$ws.Range("L114").Value = 5771.357400000001
$ws.Range("M114").Value = -363.25
$ws.Range("N114").Value = -12279.3574
$ws.Range("H117").Value = 2475.3333
$ws.Range("J117").Value = 2578.6428
$ws.Range("L117").Value = 7735.928400000001
$ws.Range("N117").Value = -14619.9284
$ws.Range("H135").Value = 3678.3438
$ws.Range("I135").Value = 12850.875
$ws.Range("J135").Value = 620.8333
$ws.Range("K135").Value = 115657.875
$ws.Range("L135").Value = 5587.4997
$ws.Range("M135").Value = -113122.875
$ws.Range("N135").Value = -10657.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1248.4166
$ws.Range("I102").Value = 995.2857
$ws.Range("J102").Value = 1602.8
$ws.Range("K102").Value = 995.2857
$ws.Range("L102").Value = 1602.8
$ws.Range("M102").Value = 626.7143
$ws.Range("N102").Value = -4846.8
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2793.7144
$ws.Range("I7").Value = 1489.1111
$ws.Range("J7").Value = 5142
$ws.Range("K7").Value = 1489.1111
$ws.Range("L7").Value = 5142
$ws.Range("M7").Value = -1377.1111
$ws.Range("N7").Value = -5366
$ws.Range("H126").Value = 2793.7144
$ws.Range("I126").Value = 1489.1111
$ws.Range("J126").Value = 5142
$ws.Range("K126").Value = 4467.3333
$ws.Range("L126").Value = 15426
$ws.Range("M126").Value = -1997.3333
$ws.Range("N126").Value = -20366
$ws.Range("H132").Value = 4490.4287
$ws.Range("I132").Value = 3754.7273
$ws.Range("J132").Value = 5299.7
$ws.Range("K132").Value = 11264.1819
$ws.Range("L132").Value = 15899.1
$ws.Range("M132").Value = -8734.1819
$ws.Range("N132").Value = -20959.1
$ws.Range("H134").Value = 56482
$ws.Range("J134").Value = 56482
$ws.Range("L134").Value = 56482
$ws.Range("N134").Value = -66622

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 981
$ws.Range("I113").Value = 671.5833
$ws.Range("J113").Value = 1511.4286
$ws.Range("K113").Value = 2014.7499
$ws.Range("L113").Value = 4534.2858
$ws.Range("M113").Value = 155.2501
$ws.Range("N113").Value = -8874.2858
$ws.Range("H132").Value = 1404677.9
$ws.Range("I132").Value = 3624927
$ws.Range("J132").Value = 2415.3157
$ws.Range("K132").Value = 10874781
$ws.Range("L132").Value = 7245.9471
$ws.Range("M132").Value = -10872251
$ws.Range("N132").Value = -12305.9471
